$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 78.30768999999999
$ws.Range("I8").Value = 78.30768999999999
$ws.Range("K8").Value = 234.92307
$ws.Range("M8").Value = -95.92307
$ws.Range("H16").Value = 3500
$ws.Range("J16").Value = 3500
$ws.Range("L16").Value = 3500
$ws.Range("N16").Value = -3960
$ws.Range("H33").Value = 423.5
$ws.Range("I33").Value = 423.5
$ws.Range("K33").Value = 423.5
$ws.Range("M33").Value = -194.5
$ws.Range("H55").Value = 308.66666
$ws.Range("I55").Value = 233.66667
$ws.Range("J55").Value = 383.66666
$ws.Range("K55").Value = 233.66667
$ws.Range("L55").Value = 383.66666
$ws.Range("M55").Value = -19.66667000000001
$ws.Range("N55").Value = -811.66666
$ws.Range("H100").Value = 2714
$ws.Range("I100").Value = 3068.5
$ws.Range("J100").Value = 1532.3334
$ws.Range("K100").Value = 3068.5
$ws.Range("L100").Value = 1532.3334
$ws.Range("M100").Value = -2527.5
$ws.Range("N100").Value = -2614.3334
$ws.Range("H103").Value = 1958.4
$ws.Range("I103").Value = 750
$ws.Range("J103").Value = 2764
$ws.Range("K103").Value = 2250
$ws.Range("L103").Value = 8292
$ws.Range("M103").Value = -1664
$ws.Range("N103").Value = -9464
$ws.Range("H129").Value = 7272.857
$ws.Range("I129").Value = 5638.6
$ws.Range("K129").Value = 16915.8
$ws.Range("M129").Value = -11915.8
$ws.Range("H132").Value = 3200.375
$ws.Range("I132").Value = 2243.6191
$ws.Range("J132").Value = 5026.909
$ws.Range("K132").Value = 6730.8573
$ws.Range("L132").Value = 15080.727
$ws.Range("M132").Value = -4200.8573
$ws.Range("N132").Value = -20140.727
$ws.Range("H141").Value = 2135.4546
$ws.Range("I141").Value = 2135.4546
$ws.Range("K141").Value = 6406.3638
$ws.Range("M141").Value = -1226.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4215.8237
$ws.Range("I32").Value = 2774.3333
$ws.Range("K32").Value = 2774.3333
$ws.Range("M32").Value = -2487.3333
$ws.Range("H97").Value = 1812.5714
$ws.Range("I97").Value = 1281.3334
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 1281.3334
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -785.3334
$ws.Range("N97").Value = -5992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8577.777
$ws.Range("I86").Value = 6860
$ws.Range("K86").Value = 6860
$ws.Range("M86").Value = -5737
$ws.Range("H89").Value = 8577.777
$ws.Range("I89").Value = 6860
$ws.Range("K89").Value = 34300
$ws.Range("M89").Value = -28684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1829.4
$ws.Range("I7").Value = 987.3077
$ws.Range("K7").Value = 987.3077
$ws.Range("M7").Value = -874.3077
$ws.Range("H22").Value = 2185.4285
$ws.Range("J22").Value = 3214.75
$ws.Range("L22").Value = 3214.75
$ws.Range("N22").Value = -3914.75
$ws.Range("H31").Value = 2213.1538
$ws.Range("I31").Value = 2427.4
$ws.Range("K31").Value = 2427.4
$ws.Range("M31").Value = -2132.4
$ws.Range("H34").Value = 2213.1538
$ws.Range("I34").Value = 2427.4
$ws.Range("K34").Value = 2427.4
$ws.Range("M34").Value = -2225.4
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H105").Value = 913.3333
$ws.Range("I105").Value = 871.25
$ws.Range("J105").Value = 997.5
$ws.Range("K105").Value = 871.25
$ws.Range("L105").Value = 997.5
$ws.Range("M105").Value = 875.75
$ws.Range("N105").Value = -4491.5
$ws.Range("H134").Value = 1526.2858
$ws.Range("I134").Value = 1526.2858
$ws.Range("K134").Value = 4578.857400000001
$ws.Range("M134").Value = -2043.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 864.1111
$ws.Range("J34").Value = 994.3333
$ws.Range("L34").Value = 2982.9999
$ws.Range("N34").Value = -3150.9999
$ws.Range("H80").Value = 14999.667
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 14999.667
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 44999.001
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -46871.001
$ws.Range("H83").Value = 14999.667
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 14999.667
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 134997.003
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -144357.003
$ws.Range("H109").Value = 2666.25
$ws.Range("I109").Value = 221.66667
$ws.Range("K109").Value = 665.00001
$ws.Range("M109").Value = 374.99999
$ws.Range("H128").Value = 286719.8
$ws.Range("I128").Value = 286719.8
$ws.Range("K128").Value = 860159.3999999999
$ws.Range("M128").Value = -855179.3999999999
$ws.Range("H131").Value = 1139.9333
$ws.Range("I131").Value = 1100
$ws.Range("K131").Value = 3300
$ws.Range("M131").Value = 1740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1330.6666
$ws.Range("I102").Value = 1099.6666
$ws.Range("K102").Value = 1099.6666
$ws.Range("M102").Value = 522.3334
$ws.Range("H132").Value = 2418.5454
$ws.Range("I132").Value = 1859
$ws.Range("K132").Value = 5577
$ws.Range("M132").Value = -3047

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 945.6667
$ws.Range("I82").Value = 945.6667
$ws.Range("K82").Value = 945.6667
$ws.Range("M82").Value = -584.6667
$ws.Range("H85").Value = 945.6667
$ws.Range("I85").Value = 945.6667
$ws.Range("K85").Value = 945.6667
$ws.Range("M85").Value = 302.3333
$ws.Range("H93").Value = 1637.5
$ws.Range("I93").Value = 1637.5
$ws.Range("K93").Value = 1637.5
$ws.Range("M93").Value = -389.5
$ws.Range("H100").Value = 5874.25
$ws.Range("I100").Value = 4499
$ws.Range("K100").Value = 4499
$ws.Range("M100").Value = -3958
$ws.Range("H132").Value = 5180.7896
$ws.Range("I132").Value = 5433.5
$ws.Range("K132").Value = 16300.5
$ws.Range("M132").Value = -13770.5
$ws.Range("H136").Value = 4574.5835
$ws.Range("J136").Value = 5400
$ws.Range("L136").Value = 16200
$ws.Range("N136").Value = -21300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 703.3333
$ws.Range("I23").Value = 110
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 110
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 119
$ws.Range("N23").Value = -1458
$ws.Range("H132").Value = 841.6316
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2071
$ws.Range("I136").Value = 1412.3334
$ws.Range("K136").Value = 4237.0002
$ws.Range("M136").Value = -1687.0002
